$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "npm start"
$ws.Range("C21").Value = "2017-11-11 19:29:43"
$ws.Range("D21").Value = "2017-11-11 19:29:47"
$ws.Range("E21").Value = 556
$ws.Range("F21").Value = 432
$ws.Range("G21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 16
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0

# Row 22
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "npm start"
$ws.Range("C22").Value = "2017-11-12 08:31:45"
$ws.Range("D22").Value = "2017-11-12 08:31:49"
$ws.Range("G22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 18
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0

# Row 23
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "npm start"
$ws.Range("C23").Value = "2017-11-12 08:32:36"
$ws.Range("D23").Value = "2017-11-12 08:32:40"
$ws.Range("E23").Value = 581
$ws.Range("F23").Value = 366
$ws.Range("G23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 19
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0

# Row 24
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "npm start"
$ws.Range("C24").Value = "2017-11-12 08:33:09"
$ws.Range("D24").Value = "2017-11-12 08:33:11"
$ws.Range("E24").Value = 580
$ws.Range("F24").Value = 489
$ws.Range("G24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 20
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0

# Row 25
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "npm start"
$ws.Range("C25").Value = "2017-11-12 08:35:59"
$ws.Range("D25").Value = "2017-11-12 08:36:04"
$ws.Range("E25").Value = 533
$ws.Range("F25").Value = 421
$ws.Range("G25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 21
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0

# Row 26
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "npm start"
$ws.Range("C26").Value = "2017-11-12 08:36:16"
$ws.Range("D26").Value = "2017-11-12 08:36:25"
$ws.Range("E26").Value = 593
$ws.Range("F26").Value = 482
$ws.Range("G26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 22
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0

# Row 27
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "npm start"
$ws.Range("C27").Value = "2017-11-12 08:42:48"
$ws.Range("D27").Value = "2017-11-12 08:42:59"
$ws.Range("E27").Value = 557
$ws.Range("F27").Value = 338
$ws.Range("G27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 23
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0

# Row 28
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "npm start"
$ws.Range("C28").Value = "2017-11-12 08:43:22"
$ws.Range("D28").Value = "2017-11-12 08:43:27"
$ws.Range("E28").Value = 492
$ws.Range("F28").Value = 450
$ws.Range("G28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 24
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0

# Row 29
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = "RemoteJob"
$ws.Range("C29").Value = "2017-11-12 08:43:27"
$ws.Range("D29").Value = "2017-11-12 08:43:27"
$ws.Range("E29").Value = 492
$ws.Range("F29").Value = 450
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 1
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 25
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0

# Row 30
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "npm start"
$ws.Range("C30").Value = "2017-11-12 08:43:46"
$ws.Range("D30").Value = "2017-11-12 08:43:49"
$ws.Range("E30").Value = 558
$ws.Range("F30").Value = 513
$ws.Range("G30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 25
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0

# Row 31
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "npm start"
$ws.Range("C31").Value = "2017-11-12 08:45:07"
$ws.Range("D31").Value = "2017-11-12 08:45:09"
$ws.Range("E31").Value = 573
$ws.Range("F31").Value = 502
$ws.Range("G31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 26
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0

# Row 32
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "npm start"
$ws.Range("C32").Value = "2017-11-12 08:49:37"
$ws.Range("D32").Value = "2017-11-12 08:49:41"
$ws.Range("E32").Value = 550
$ws.Range("F32").Value = 495
$ws.Range("G32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 27
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0

# Row 33
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "npm start"
$ws.Range("C33").Value = "2017-11-12 08:53:14"
$ws.Range("D33").Value = "2017-11-12 08:53:20"
$ws.Range("E33").Value = 542
$ws.Range("F33").Value = 498
$ws.Range("G33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 28
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0

# Row 34
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = "npm start"
$ws.Range("C34").Value = "2017-11-12 08:53:30"
$ws.Range("D34").Value = "2017-11-12 08:53:33"
$ws.Range("E34").Value = 584
$ws.Range("F34").Value = 492
$ws.Range("G34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 29
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0

# Row 35
$ws.Range("A35").Value = 34
$ws.Range("B35").Value = "npm start"
$ws.Range("C35").Value = "2017-11-12 08:53:47"
$ws.Range("D35").Value = "2017-11-12 08:53:54"
$ws.Range("E35").Value = 498
$ws.Range("F35").Value = 366
$ws.Range("G35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 30
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 0

# Row 36
$ws.Range("A36").Value = 35
$ws.Range("B36").Value = "RemoteJob"
$ws.Range("C36").Value = "2017-11-12 08:53:54"
$ws.Range("D36").Value = "2017-11-12 08:53:58"
$ws.Range("E36").Value = 498
$ws.Range("F36").Value = 366
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = 1
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 31
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 0

# Row 37
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = "npm start"
$ws.Range("C37").Value = "2017-11-12 08:54:14"
$ws.Range("D37").Value = "2017-11-12 08:54:16"
$ws.Range("E37").Value = 415
$ws.Range("F37").Value = 389
$ws.Range("G37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 31
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 0

# Row 38
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = "npm start"
$ws.Range("C38").Value = "2017-11-12 08:57:30"
$ws.Range("D38").Value = "2017-11-12 08:57:35"
$ws.Range("E38").Value = 492
$ws.Range("F38").Value = 419
$ws.Range("G38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 32
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 0

# Row 39
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "npm start"
$ws.Range("C39").Value = "2017-11-12 08:59:39"
$ws.Range("D39").Value = "2017-11-12 08:59:47"
$ws.Range("E39").Value = 508
$ws.Range("F39").Value = 445
$ws.Range("G39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 33
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 0

# Row 40
$ws.Range("A40").Value = 39
$ws.Range("B40").Value = "npm start"
$ws.Range("C40").Value = "2017-11-12 09:00:25"
$ws.Range("D40").Value = "2017-11-12 09:00:29"
$ws.Range("E40").Value = 562
$ws.Range("F40").Value = 504
$ws.Range("G40").Value = 0
$ws.Range("I40").Value = 4
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 34
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = 0

# Row 41
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = "npm start"
$ws.Range("C41").Value = "2017-11-12 09:02:59"
$ws.Range("D41").Value = "2017-11-12 09:03:13"
$ws.Range("E41").Value = 499
$ws.Range("F41").Value = 260
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = 1
$ws.Range("I41").Value = 2
$ws.Range("J41").Value = 0.8
$ws.Range("K41").Value = 35
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = 0

# Row 42
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = "generateImageFromDB.js — ~/DevHouse/office/RemoteJob/desktop-app — Atom"
$ws.Range("C42").Value = "2017-11-12 09:03:13"
$ws.Range("D42").Value = "2017-11-12 09:03:19"
$ws.Range("E42").Value = 164
$ws.Range("F42").Value = 250
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = 2
$ws.Range("I42").Value = 4
$ws.Range("J42").Value = 1.2
$ws.Range("K42").Value = 36
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0

# Row 43
$ws.Range("A43").Value = 42
$ws.Range("B43").Value = "Project — ~/DevHouse/office/RemoteJob/desktop-app — Atom"
$ws.Range("C43").Value = "2017-11-12 09:03:19"
$ws.Range("D43").Value = "2017-11-12 09:03:34"
$ws.Range("E43").Value = 516
$ws.Range("F43").Value = 487
$ws.Range("G43").Value = 3
$ws.Range("H43").Value = 1
$ws.Range("I43").Value = 4
$ws.Range("J43").Value = 2
$ws.Range("K43").Value = 36
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = 0

# Row 44
$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "RemoteJob"
$ws.Range("C44").Value = "2017-11-12 09:03:34"
$ws.Range("D44").Value = "2017-11-12 09:03:37"
$ws.Range("E44").Value = 140
$ws.Range("F44").Value = 360
$ws.Range("G44").Value = 5
$ws.Range("H44").Value = 1
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 1.65
$ws.Range("K44").Value = 36
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0

# Row 45
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "x-caja-desktop"
$ws.Range("C45").Value = "2017-11-12 09:03:37"
$ws.Range("D45").Value = "2017-11-12 09:03:38"
$ws.Range("E45").Value = 569
$ws.Range("F45").Value = 473
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 1
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0.33
$ws.Range("K45").Value = 36
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = 0

Write-Output "done"